$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Soliera report update ("aggiornato a 2/3, aggiornati i report")
#
# 1) A new daily data point for date 44235 (2021-02-03) with 2 new cases was
#    inserted between the existing rows for 44234 and 44236. This shifts all
#    subsequent rows down by one.
# 2) The 7-day rolling window values (col C = somma mobile 7gg., col D = per
#    100k abitanti) shift along with it: rows whose window now includes the
#    new data point get recomputed totals, and the trailing "not yet
#    computable" rows move down by one row as well (one more row near the
#    bottom now has a computed total that previously did not).
# 3) A brand-new trailing row for date 44257 (2021-03-07) with 5 new cases
#    was appended at the end (its rolling-sum columns are not yet computed,
#    matching the existing trailing rows).
# ---------------------------------------------------------------------------

# 1) Insert a new row at 93 -- shifts old rows 93:113 down to 94:114, values
#    and formatting intact.
$ws.Rows("93:93").Insert()

# Fix up the style of the newly-inserted A93 cell: Insert() stamps the new
# row with a freshly-minted style, but this column should reuse the same
# style as the surrounding date cells (style index 2 in the original file).
$ws.Range("A92").Copy()
$ws.Range("A93").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the newly-inserted row for date 44235.
$ws.Range("A93").Value2 = 44235
$ws.Range("B93").Value2 = 2
$ws.Range("C93").Value2 = 13
$ws.Range("D93").Value2 = 83.96305625524769

# 2) Recompute the 7-day rolling window figures for the rows whose window
#    now includes the newly-inserted date (44232-44237).
$ws.Range("C90").Value2 = 11
$ws.Range("D90").Value2 = 71.04566298520959

$ws.Range("C91").Value2 = 12
$ws.Range("D91").Value2 = 77.50435962022863

$ws.Range("C92").Value2 = 12
$ws.Range("D92").Value2 = 77.50435962022863

$ws.Range("C94").Value2 = 16
$ws.Range("D94").Value2 = 103.3391461603049

$ws.Range("C95").Value2 = 20
$ws.Range("D95").Value2 = 129.1739327003811

# The trailing "rolling sum not available yet" rows also shift down by one:
# the row for date 44254 (now row 112) gains a computed total it did not
# have before.
$ws.Range("C112").Value2 = 43
$ws.Range("D112").Value2 = 277.7239553058193

# 3) Append a brand-new trailing row for date 44257 at the end of the table
#    (row 115), mirroring the still-unresolved rolling sum of the other
#    recent trailing rows (left blank, like rows 113/114 above it).
$ws.Range("A114").Copy()
$ws.Range("A115").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A115").Value2 = 44257
$ws.Range("B115").Value2 = 5
$ws.Range("C115").Value2 = ""
$ws.Range("D115").Value2 = ""

Write-Host "Soliera report updated through 2021-03-07 (dimension now A1:D115)."
